$wb = $excel.ActiveWorkbook

# Duplicate the "nativity" sheet (to inherit its sheetFormatPr / namespace
# cosmetics), drop it right after "lep", rename it to "pop", and wipe its
# old contents + styles before writing the new population data.
$src = $wb.Worksheets.Item("nativity")
$afterTarget = $wb.Worksheets.Item("lep")
[void]$src.Copy($null, $afterTarget)
$newSheet = $wb.Worksheets.Item("nativity (2)")
$newSheet.Name = "pop"
[void]$newSheet.Cells.Clear()

$header = New-Object 'object[,]' 1,8
$header[0,0] = "geoid"
$header[0,1] = "geodisplaylabel"
$header[0,2] = "popgroupid"
$header[0,3] = "est_tot_pop"
$header[0,4] = "moe_tot_pop"
$header[0,5] = "checking_tot_pop"
$header[0,6] = "est_national_pop"
$header[0,7] = "pct_pop"
$newSheet.Range("A1:H1").Value = $header

$popData = New-Object 'object[,]' 71,8
$popData[0,0] = "0100000US"
$popData[0,1] = "United States"
$popData[0,2] = 1
$popData[0,3] = 316515008
$popData[0,4] = 0
$popData[0,5] = 0
$popData[0,6] = 316515008
$popData[0,7] = 1
$popData[1,0] = "0100000US"
$popData[1,1] = "United States"
$popData[1,2] = 4
$popData[1,3] = 39908096
$popData[1,4] = 26702
$popData[1,5] = 0
$popData[1,6] = 316515008
$popData[1,7] = 0.12608595192432404
$popData[2,0] = "0100000US"
$popData[2,1] = "United States"
$popData[2,2] = 6
$popData[2,3] = 2569170
$popData[2,4] = 14914
$popData[2,5] = 0
$popData[2,6] = 316515008
$popData[2,7] = 0.008117055520415306
$popData[3,0] = "0100000US"
$popData[3,1] = "United States"
$popData[3,2] = 12
$popData[3,3] = 16235305
$popData[3,4] = 20725
$popData[3,5] = 0
$popData[3,6] = 316515008
$popData[3,7] = 0.05129395052790642
$popData[4,0] = "0100000US"
$popData[4,1] = "United States"
$popData[4,2] = 13
$popData[4,3] = 3303512
$popData[4,4] = 17090
$popData[4,5] = 0
$popData[4,6] = 316515008
$popData[4,7] = 0.010437142103910446
$popData[5,0] = "0100000US"
$popData[5,1] = "United States"
$popData[5,2] = 14
$popData[5,3] = 153234
$popData[5,4] = 5479
$popData[5,5] = 0
$popData[5,6] = 316515008
$popData[5,7] = 0.00048412871547043324
$popData[6,0] = "0100000US"
$popData[6,1] = "United States"
$popData[6,2] = 15
$popData[6,3] = 263396
$popData[6,4] = 8409
$popData[6,5] = 0
$popData[6,6] = 316515008
$popData[6,7] = 0.0008321753703057766
$popData[7,0] = "0100000US"
$popData[7,1] = "United States"
$popData[7,2] = 16
$popData[7,3] = 3852099
$popData[7,4] = 19831
$popData[7,5] = 0
$popData[7,6] = 316515008
$popData[7,7] = 0.012170352041721344
$popData[8,0] = "0100000US"
$popData[8,1] = "United States"
$popData[8,2] = 17
$popData[8,3] = 3693977
$popData[8,4] = 18855
$popData[8,5] = 0
$popData[8,6] = 316515008
$popData[8,7] = 0.01167078036814928
$popData[9,0] = "0100000US"
$popData[9,1] = "United States"
$popData[9,2] = 18
$popData[9,3] = 148275
$popData[9,4] = 4247
$popData[9,5] = 0
$popData[9,6] = 316515008
$popData[9,7] = 0.00046846119221299887
$popData[10,0] = "0100000US"
$popData[10,1] = "United States"
$popData[10,2] = 19
$popData[10,3] = 2717844
$popData[10,4] = 17934
$popData[10,5] = 0
$popData[10,6] = 316515008
$popData[10,7] = 0.008586777374148369
$popData[11,0] = "0100000US"
$popData[11,1] = "United States"
$popData[11,2] = 20
$popData[11,3] = 267009
$popData[11,4] = 6734
$popData[11,5] = 0
$popData[11,6] = 316515008
$popData[11,7] = 0.0008435903582721949
$popData[12,0] = "0100000US"
$popData[12,1] = "United States"
$popData[12,2] = 21
$popData[12,3] = 71451
$popData[12,4] = 2920
$popData[12,5] = 0
$popData[12,6] = 316515008
$popData[12,7] = 0.00022574285685550421
$popData[13,0] = "0100000US"
$popData[13,1] = "United States"
$popData[13,2] = 22
$popData[13,3] = 779637
$popData[13,4] = 9694
$popData[13,5] = 0
$popData[13,6] = 316515008
$popData[13,7] = 0.002463191282004118
$popData[14,0] = "0100000US"
$popData[14,1] = "United States"
$popData[14,2] = 23
$popData[14,3] = 1460214
$popData[14,4] = 13339
$popData[14,5] = 0
$popData[14,6] = 316515008
$popData[14,7] = 0.004613411612808704
$popData[15,0] = "0100000US"
$popData[15,1] = "United States"
$popData[15,2] = 24
$popData[15,3] = 207999
$popData[15,4] = 4997
$popData[15,5] = 0
$popData[15,6] = 316515008
$popData[15,7] = 0.0006571536650881171
$popData[16,0] = "0100000US"
$popData[16,1] = "United States"
$popData[16,2] = 25
$popData[16,3] = 18803
$popData[16,4] = 1544
$popData[16,5] = 0
$popData[16,6] = 316515008
$popData[16,7] = 0.00005940634582657367
$popData[17,0] = "0100000US"
$popData[17,1] = "United States"
$popData[17,2] = 26
$popData[17,3] = 414880
$popData[17,4] = 11430
$popData[17,5] = 0
$popData[17,6] = 316515008
$popData[17,7] = 0.0013107751728966832
$popData[18,0] = "0100000US"
$popData[18,1] = "United States"
$popData[18,2] = 27
$popData[18,3] = 46036
$popData[18,4] = 2724
$popData[18,5] = 0
$popData[18,6] = 316515008
$popData[18,7] = 0.00014544649457093328
$popData[19,0] = "0100000US"
$popData[19,1] = "United States"
$popData[19,2] = 28
$popData[19,3] = 188673
$popData[19,4] = 4784
$popData[19,5] = 0
$popData[19,6] = 316515008
$popData[19,7] = 0.0005960949347354472
$popData[20,0] = "0100000US"
$popData[20,1] = "United States"
$popData[20,2] = 29
$popData[20,3] = 1710547
$popData[20,4] = 18780
$popData[20,5] = 0
$popData[20,6] = 316515008
$popData[20,7] = 0.005404315423220396
$popData[21,0] = "0100000US"
$popData[21,1] = "United States"
$popData[21,2] = 31
$popData[21,3] = 19167716
$popData[21,4] = 19485
$popData[21,5] = 0
$popData[21,6] = 316515008
$popData[21,7] = 0.060558632016181946
$popData[22,0] = "0100000US"
$popData[22,1] = "United States"
$popData[22,2] = 32
$popData[22,3] = 3590279
$popData[22,4] = 16777
$popData[22,5] = 0
$popData[22,6] = 316515008
$popData[22,7] = 0.011343155987560749
$popData[23,0] = "0100000US"
$popData[23,1] = "United States"
$popData[23,2] = 33
$popData[23,3] = 164821
$popData[23,4] = 5618
$popData[23,5] = 0
$popData[23,6] = 316515008
$popData[23,7] = 0.0005207367357797921
$popData[24,0] = "0100000US"
$popData[24,1] = "United States"
$popData[24,2] = 34
$popData[24,3] = 316640
$popData[24,4] = 8596
$popData[24,5] = 0
$popData[24,6] = 316515008
$popData[24,7] = 0.0010003949282690883
$popData[25,0] = "0100000US"
$popData[25,1] = "United States"
$popData[25,2] = 35
$popData[25,3] = 4597905
$popData[25,4] = 24921
$popData[25,5] = 0
$popData[25,6] = 316515008
$popData[25,7] = 0.014526656828820705
$popData[26,0] = "0100000US"
$popData[26,1] = "United States"
$popData[26,2] = 36
$popData[26,3] = 4428363
$popData[26,4] = 23401
$popData[26,5] = 0
$popData[26,6] = 316515008
$popData[26,7] = 0.013991004787385464
$popData[27,0] = "0100000US"
$popData[27,1] = "United States"
$popData[27,2] = 37
$popData[27,3] = 181029
$popData[27,4] = 4898
$popData[27,5] = 0
$popData[27,6] = 316515008
$popData[27,7] = 0.0005719444598071277
$popData[28,0] = "0100000US"
$popData[28,1] = "United States"
$popData[28,2] = 38
$popData[28,3] = 3707082
$popData[28,4] = 24374
$popData[28,5] = 0
$popData[28,6] = 316515008
$popData[28,7] = 0.011712184175848961
$popData[29,0] = "0100000US"
$popData[29,1] = "United States"
$popData[29,2] = 39
$popData[29,3] = 280410
$popData[29,4] = 6971
$popData[29,5] = 0
$popData[29,6] = 316515008
$popData[29,7] = 0.0008859295630827546
$popData[30,0] = "0100000US"
$popData[30,1] = "United States"
$popData[30,2] = 40
$popData[30,3] = 108477
$popData[30,4] = 3749
$popData[30,5] = 0
$popData[30,6] = 316515008
$popData[30,7] = 0.00034272309858351946
$popData[31,0] = "0100000US"
$popData[31,1] = "United States"
$popData[31,2] = 41
$popData[31,3] = 1388163
$popData[31,4] = 15805
$popData[31,5] = 0
$popData[31,6] = 316515008
$popData[31,7] = 0.004385773092508316
$popData[32,0] = "0100000US"
$popData[32,1] = "United States"
$popData[32,2] = 42
$popData[32,3] = 1792437
$popData[32,4] = 14297
$popData[32,5] = 0
$popData[32,6] = 316515008
$popData[32,7] = 0.0056630396284163
$popData[33,0] = "0100000US"
$popData[33,1] = "United States"
$popData[33,2] = 43
$popData[33,3] = 256416
$popData[33,4] = 5432
$popData[33,5] = 0
$popData[33,6] = 316515008
$popData[33,7] = 0.0008101226994767785
$popData[34,0] = "0100000US"
$popData[34,1] = "United States"
$popData[34,2] = 44
$popData[34,3] = 28818
$popData[34,4] = 1890
$popData[34,5] = 0
$popData[34,6] = 316515008
$popData[34,7] = 0.00009104781202040613
$popData[35,0] = "0100000US"
$popData[35,1] = "United States"
$popData[35,2] = 45
$popData[35,3] = 455005
$popData[35,4] = 11349
$popData[35,5] = 0
$popData[35,6] = 316515008
$popData[35,7] = 0.0014375463360920548
$popData[36,0] = "0100000US"
$popData[36,1] = "United States"
$popData[36,2] = 46
$popData[36,3] = 52267
$popData[36,4] = 2925
$popData[36,5] = 0
$popData[36,6] = 316515008
$popData[36,7] = 0.00016513277660124004
$popData[37,0] = "0100000US"
$popData[37,1] = "United States"
$popData[37,2] = 47
$popData[37,3] = 278281
$popData[37,4] = 5520
$popData[37,5] = 0
$popData[37,6] = 316515008
$popData[37,7] = 0.0008792032022029161
$popData[38,0] = "0100000US"
$popData[38,1] = "United States"
$popData[38,2] = 48
$popData[38,3] = 1928363
$popData[38,4] = 19234
$popData[38,5] = 0
$popData[38,6] = 316515008
$popData[38,7] = 0.006092485040426254
$popData[39,0] = "0100000US"
$popData[39,1] = "United States"
$popData[39,2] = 50
$popData[39,3] = 546255
$popData[39,4] = 4552
$popData[39,5] = 0
$popData[39,6] = 316515008
$popData[39,7] = 0.0017258423613384366
$popData[40,0] = "0100000US"
$popData[40,1] = "United States"
$popData[40,2] = 51
$popData[40,3] = 344487
$popData[40,4] = 4839
$popData[40,5] = 0
$popData[40,6] = 316515008
$popData[40,7] = 0.0010883748764172196
$popData[41,0] = "0100000US"
$popData[41,1] = "United States"
$popData[41,2] = 52
$popData[41,3] = 174460
$popData[41,4] = 4519
$popData[41,5] = 0
$popData[41,6] = 316515008
$popData[41,7] = 0.0005511902854777873
$popData[42,0] = "0100000US"
$popData[42,1] = "United States"
$popData[42,2] = 53
$popData[42,3] = 109455
$popData[42,4] = 3835
$popData[42,5] = 0
$popData[42,6] = 316515008
$popData[42,7] = 0.00034581299405544996
$popData[43,0] = "0100000US"
$popData[43,1] = "United States"
$popData[43,2] = 54
$popData[43,3] = 45453
$popData[43,4] = 2981
$popData[43,5] = 0
$popData[43,6] = 316515008
$popData[43,7] = 0.00014360457134898752
$popData[44,0] = "0100000US"
$popData[44,1] = "United States"
$popData[44,2] = 55
$popData[44,3] = 138360
$popData[44,4] = 3689
$popData[44,5] = 0
$popData[44,6] = 316515008
$popData[44,7] = 0.0004371356626506895
$popData[45,0] = "0100000US"
$popData[45,1] = "United States"
$popData[45,2] = 56
$popData[45,3] = 73088
$popData[45,4] = 2797
$popData[45,5] = 0
$popData[45,6] = 316515008
$popData[45,7] = 0.00023091479670256376
$popData[46,0] = "0100000US"
$popData[46,1] = "United States"
$popData[46,2] = 57
$popData[46,3] = 34239
$popData[46,4] = 1993
$popData[46,5] = 0
$popData[46,6] = 316515008
$popData[46,7] = 0.00010817496513482183
$popData[47,0] = "0100000US"
$popData[47,1] = "United States"
$popData[47,2] = 58
$popData[47,3] = 33468
$popData[47,4] = 1972
$popData[47,5] = 0
$popData[47,6] = 316515008
$popData[47,7] = 0.00010573906183708459
$popData[48,0] = "0100000US"
$popData[48,1] = "United States"
$popData[48,2] = 60
$popData[48,3] = 1262434
$popData[48,4] = 9700
$popData[48,5] = 0
$popData[48,6] = 316515008
$popData[48,7] = 0.00398854399099946
$popData[49,0] = "0100000US"
$popData[49,1] = "United States"
$popData[49,2] = 61
$popData[49,3] = 783326
$popData[49,4] = 7279
$popData[49,5] = 0
$popData[49,6] = 316515008
$popData[49,7] = 0.0024748463183641434
$popData[50,0] = "0100000US"
$popData[50,1] = "United States"
$popData[50,2] = 62
$popData[50,3] = 549858
$popData[50,4] = 8077
$popData[50,5] = 0
$popData[50,6] = 316515008
$popData[50,7] = 0.001737225684337318
$popData[51,0] = "0100000US"
$popData[51,1] = "United States"
$popData[51,2] = 63
$popData[51,3] = 182968
$popData[51,4] = 5079
$popData[51,5] = 0
$popData[51,6] = 316515008
$popData[51,7] = 0.0005780705250799656
$popData[52,0] = "0100000US"
$popData[52,1] = "United States"
$popData[52,2] = 64
$popData[52,3] = 62458
$popData[52,4] = 3119
$popData[52,5] = 0
$popData[52,6] = 316515008
$popData[52,7] = 0.00019733030057977885
$popData[53,0] = "0100000US"
$popData[53,1] = "United States"
$popData[53,2] = 65
$popData[53,3] = 207128
$popData[53,4] = 4178
$popData[53,5] = 0
$popData[53,6] = 316515008
$popData[53,7] = 0.000654401839710772
$popData[54,0] = "0100000US"
$popData[54,1] = "United States"
$popData[54,2] = 66
$popData[54,3] = 130476
$popData[54,4] = 3970
$popData[54,5] = 0
$popData[54,6] = 316515008
$popData[54,7] = 0.00041222688741981983
$popData[55,0] = "0100000US"
$popData[55,1] = "United States"
$popData[55,2] = 67
$popData[55,3] = 43211
$popData[55,4] = 2222
$popData[55,5] = 0
$popData[55,6] = 316515008
$popData[55,7] = 0.00013652117922902107
$popData[56,0] = "0100000US"
$popData[56,1] = "United States"
$popData[56,2] = 68
$popData[56,3] = 42110
$popData[56,4] = 2220
$popData[56,5] = 0
$popData[56,6] = 316515008
$popData[56,7] = 0.0001330426603090018
$popData[57,0] = "0100000US"
$popData[57,1] = "United States"
$popData[57,2] = 72
$popData[57,3] = 21441
$popData[57,4] = 2416
$popData[57,5] = 0
$popData[57,6] = 316515008
$popData[57,7] = 0.0000677408606861718
$popData[58,0] = "0100000US"
$popData[58,1] = "United States"
$popData[58,2] = 73
$popData[58,3] = 126590
$popData[58,4] = 5388
$popData[58,5] = 0
$popData[58,6] = 316515008
$popData[58,7] = 0.0003999494365416467
$popData[59,0] = "0100000US"
$popData[59,1] = "United States"
$popData[59,2] = 75
$popData[59,3] = 15919
$popData[59,4] = 1626
$popData[59,5] = 0
$popData[59,6] = 316515008
$popData[59,7] = 0.000050294613174628466
$popData[60,0] = "0100000US"
$popData[60,1] = "United States"
$popData[60,2] = 76
$popData[60,3] = 103526
$popData[60,4] = 4285
$popData[60,5] = 0
$popData[60,6] = 316515008
$popData[60,7] = 0.000327080866554752
$popData[61,0] = "0100000US"
$popData[61,1] = "United States"
$popData[61,2] = 80
$popData[61,3] = 24875
$popData[61,4] = 2425
$popData[61,5] = 0
$popData[61,6] = 316515008
$popData[61,7] = 0.00007859026663936675
$popData[62,0] = "0100000US"
$popData[62,1] = "United States"
$popData[62,2] = 81
$popData[62,3] = 138426
$popData[62,4] = 5576
$popData[62,5] = 0
$popData[62,6] = 316515008
$popData[62,7] = 0.000437344191595912
$popData[63,0] = "0100000US"
$popData[63,1] = "United States"
$popData[63,2] = 83
$popData[63,3] = 20546
$popData[63,4] = 1898
$popData[63,5] = 0
$popData[63,6] = 316515008
$popData[63,7] = 0.00006491319072665647
$popData[64,0] = "0100000US"
$popData[64,1] = "United States"
$popData[64,2] = 84
$popData[64,3] = 111021
$popData[64,4] = 4578
$popData[64,5] = 0
$popData[64,6] = 316515008
$popData[64,7] = 0.00035076061612926424
$popData[65,0] = "0100000US"
$popData[65,1] = "United States"
$popData[65,2] = 85
$popData[65,3] = 11627
$popData[65,4] = 870
$popData[65,5] = 0
$popData[65,6] = 316515008
$popData[65,7] = 0.000036734436434926465
$popData[66,0] = "0100000US"
$popData[66,1] = "United States"
$popData[66,2] = 96
$popData[66,3] = 23444
$popData[66,4] = 1988
$popData[66,5] = 0
$popData[66,6] = 316515008
$popData[66,7] = 0.0000740691612008959
$popData[67,0] = "0100000US"
$popData[67,1] = "United States"
$popData[67,2] = 176
$popData[67,3] = 26856
$popData[67,4] = 2368
$popData[67,5] = 0
$popData[67,6] = 316515008
$popData[67,7] = 0.00008484905993100256
$popData[68,0] = "0100000US"
$popData[68,1] = "United States"
$popData[68,2] = 177
$popData[68,3] = 8957
$popData[68,4] = 1133
$popData[68,5] = 0
$popData[68,6] = 316515008
$popData[68,7] = 0.000028298816687311046
$popData[69,0] = "0100000US"
$popData[69,1] = "United States"
$popData[69,2] = 400
$popData[69,3] = 54232204
$popData[69,4] = 2036
$popData[69,5] = 0
$popData[69,6] = 316515008
$popData[69,7] = 0.1713416576385498
$popData[70,0] = "0100000US"
$popData[70,1] = "United States"
$popData[70,2] = 451
$popData[70,3] = 197258272
$popData[70,4] = 8320
$popData[70,5] = 0
$popData[70,6] = 316515008
$popData[70,7] = 0.6232193112373352
$newSheet.Range("A2:H72").Value = $popData

[void]$newSheet.Range("A1:H72").Select()
[void]$newSheet.Activate()

Write-Output "pop sheet populated"
